# Update from CTP review round 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content updates ---

# report_id Type: integer -> number
$ws.Range("C2").Value = "number"

# submission_date Description
$ws.Range("D3").Value = "Date report was received by CTP; this is the earliest date of report receipt, either to Safety Reporting Portal (SRP) or by other means"

# number_tobacco_products Type: integer -> number
$ws.Range("C7").Value = "number"

# number_product_problems Type: integer -> number
$ws.Range("C8").Value = "number"

# number_health_problems Type: integer -> number, and fix Description typo
$ws.Range("C9").Value = "number"
$ws.Range("D9").Value = "System-calculated number of Health Problems (i.e., MedDRA terms selected from a standardized list of symptoms, signs, diagnoses and outcomes) reported, displayed as a whole number, ≥0."

# nonuser_affected Description: append question-timing clarification
$ws.Range("D10").Value = "Displays text reflecting the response to this optional question (2017 - 12/14/2018) or required question (12/15/2018 onward) as “No information provided” if not answered, or Yes/No."

# --- Row height updates (auto-fit wrap-text rows to new content) ---
$ws.Rows.Item(3).RowHeight = 34
$ws.Rows.Item(4).RowHeight = 17
$ws.Rows.Item(5).RowHeight = 17
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(10).RowHeight = 51

# --- Selection update ---
$ws.Range("F10").Select()

# --- Page setup ---
$ws.PageSetup.Zoom = 65
